$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "@Alias_XX" table (A1:E5) is being replaced with a new,
# differently-shaped alias table (A1:E4). Clear the old range first so the
# row-5 leftovers don't survive.
$ws.Range("A1:E5").ClearContents()

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Fields"
$ws.Range("C1").Value = "Csharp"
$ws.Range("D1").Value = "Go"
$ws.Range("E1").Value = "Cpp"

# Row 2 - Vector3 alias
$ws.Range("A2").Value = "V3"
$ws.Range("B2").Value = "x|y|z"
$ws.Range("C2").Value = "UnityEngine.Vector3"
$ws.Range("D2").Value = "*gmath.Vector3"

# Row 3 - PairItemIntBool alias
$ws.Range("A3").Value = "PairItemIntBool"
$ws.Range("B3").Value = "Id|Exist"
$ws.Range("C3").Value = "PairItemIntBool"

# Row 4 - PairItemIntInt64 alias
$ws.Range("A4").Value = "PairItemIntInt64"
$ws.Range("B4").Value = "Id|Value"

# Column width tweaks (only columns 1, 4 and 5 actually change width;
# columns 2, 3 and 6 are left untouched so column 3's bestFit flag survives).
$ws.Columns.Item(1).ColumnWidth = 18.142857142857142
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(5).ColumnWidth = 18.142857142857142

# Match the selection recorded in the saved file
$ws.Range("B5").Select() | Out-Null
